$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = 7312
}
for ($r = 33; $r -le 151; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}
for ($r = 152; $r -le 185; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}
